$d = $word.ActiveDocument

$old = "In this stage of testing the application has been sent to some doctors to be tested on clinical and hospital and test the system efficiency and outputs is correct or not, retrieve feedback to our team."
$new = "At this stage we will test the application on real users, farmers and people who have some plants at home, retrieve feedback to our team."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Find/Replace result: $found"
